# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets contain the same event list; only the F6 starting value differs
# between them (219 vs 220) but both are updated to the same new value (221).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 254
    3  = 1310
    4  = 142
    6  = 221
    7  = 94
    8  = 11
    9  = 175
    10 = 125
    11 = 4418
    12 = 6688
    14 = 53
    15 = 94
    16 = 561
    18 = 4096
    19 = 459
    20 = 66
    21 = 46
    22 = 2670
    25 = 162
    26 = 339
    29 = 214
    30 = 30
    31 = 1607
    32 = 1012
    34 = 119
    36 = 532
    40 = 620
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
